# "fixed base spd of spectra"
#
# The "Speed Calculators" sheet computes each monster's combat speed from
# a base-speed input in column D plus two global multipliers entered in
# row 3 (speed-leader % in E3, speed-totem % in F3):
#     E{row} = D{row} * (1 + $E$3 + $F$3)
#     I{row} = (E{row} + G{row}) * (1 + H{row})
#     K{row}:U{row} = I{row} * (0.07 * tick#)
#
# This edit:
#   1) corrects Spectra's (row 5) base speed value, D5, from 116 -> 126
#   2) sets the "speed leader %" input, E3, from 0 -> 0.24
#
# Everything else on the sheet (E4:E21, I4:I21, K4:U21, ...) is a plain
# formula, so it recalculates automatically once these two inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Speed Calculators")

# Speed leader % input (row 3) - used by every row's SUM($D*(1+$E$3+$F$3))
$ws.Range("E3").Value = 0.24

# Spectra's corrected base speed (row 5, column D)
$ws.Range("D5").Value = 126

# Leave the selection where the author last clicked while making this edit
$ws.Range("E3").Select() | Out-Null
